$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row for consistency: "header" -> "subject", "filepath" -> "attachment"
$ws.Range("D1").Value = "attachment"
$ws.Range("B1").Value = "subject"

# Move the active selection to C10
$null = $ws.Range("C10").Select()
